$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = "Handback transform failed"
$wsOverview.Range("F7").Value = "Handback transform failed"

# --- Sheet "zh-cn" -------------------------------------------------------
# NOTE: ColumnWidth is expressed in "characters" and Excel round-trips it
# through a pixel-grid conversion when serializing the OOXML `width`
# attribute. 39.17 is the character-width value that serializes back to a
# stored width of exactly 40 (matching the target column width).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("P7").Value = "Handback file name: qtibu2ux.xmh is different with handoff file name: 37d3bf1a-24db-4978-81e0-5734adabaaa2.6cb22e404def2865d2ceefff17308b264bc5b26f.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- Sheet "de-de" ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("P7").Value = "Handback file name: qtibu2ux.xmh is different with handoff file name: 37d3bf1a-24db-4978-81e0-5734adabaaa2.6cb22e404def2865d2ceefff17308b264bc5b26f.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
